$wb = $excel.ActiveWorkbook

# --- ManagerDetails: update the manager record (email/name/id) ---
$managerDetails = $wb.Worksheets.Item("ManagerDetails")
$managerDetails.Range("A2").Value = "Wiley@Wiley.com"
$managerDetails.Range("B2").Value = "Wiley Edge"
$managerDetails.Range("C2").Value = "'99999"

# --- Insert a new "MappingCourseManager" sheet before "C343 19.03.23" ---
$c343_1903 = $wb.Worksheets.Item("C343 19.03.23")

# Move the selection on the existing "C343 19.03.23" sheet (it will stop
# being the active sheet once the new sheet below is inserted/activated).
$c343_1903.Range("I26").Select() | Out-Null

$mappingCourseManager = $wb.Worksheets.Add($c343_1903, $null)
$mappingCourseManager.Name = "MappingCourseManager"
$mappingCourseManager.Range("A2").Value = "C343"
$mappingCourseManager.Range("B2").Value = "Wiley@Wiley.com"
$mappingCourseManager.Range("K18").Select() | Out-Null
